$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44722
$ws.Cells.Item(2, 10).Value = 150
$ws.Cells.Item(2, 11).Value = 18000
$ws.Cells.Item(2, 12).Value = 20000
$ws.Cells.Item(2, 13).Value = 18933
$ws.Cells.Item(2, 16).Value = 1262

$ws.Cells.Item(3, 4).Value = 44736
$ws.Cells.Item(3, 10).Value = 180
$ws.Cells.Item(3, 11).Value = 17000
$ws.Cells.Item(3, 12).Value = 19000
$ws.Cells.Item(3, 13).Value = 17889
$ws.Cells.Item(3, 16).Value = 1193

$ws.Cells.Item(4, 4).Value = 44400
$ws.Cells.Item(4, 10).Value = 130
$ws.Cells.Item(4, 11).Value = 24000
$ws.Cells.Item(4, 12).Value = 24000
$ws.Cells.Item(4, 13).Value = 24000
$ws.Cells.Item(4, 16).Value = 1600

$ws.Cells.Item(5, 4).Value = 44446
$ws.Cells.Item(5, 10).Value = 150
$ws.Cells.Item(5, 11).Value = 22000
$ws.Cells.Item(5, 12).Value = 24000
$ws.Cells.Item(5, 13).Value = 22667
$ws.Cells.Item(5, 16).Value = 1511

$ws.Cells.Item(6, 4).Value = 44742
$ws.Cells.Item(6, 10).Value = 400
$ws.Cells.Item(6, 11).Value = 18000
$ws.Cells.Item(6, 12).Value = 20000
$ws.Cells.Item(6, 13).Value = 18850
$ws.Cells.Item(6, 16).Value = 1257

$ws.Cells.Item(7, 4).Value = 44392
$ws.Cells.Item(7, 10).Value = 220
$ws.Cells.Item(7, 11).Value = 23000
$ws.Cells.Item(7, 12).Value = 23000
$ws.Cells.Item(7, 13).Value = 23000
$ws.Cells.Item(7, 16).Value = 1533

$ws.Cells.Item(8, 4).Value = 44483
$ws.Cells.Item(8, 10).Value = 220
$ws.Cells.Item(8, 11).Value = 18000
$ws.Cells.Item(8, 12).Value = 20000
$ws.Cells.Item(8, 13).Value = 18909
$ws.Cells.Item(8, 16).Value = 1261

$ws.Cells.Item(9, 4).Value = 44749
$ws.Cells.Item(9, 10).Value = 220
$ws.Cells.Item(9, 11).Value = 18000
$ws.Cells.Item(9, 12).Value = 20000
$ws.Cells.Item(9, 13).Value = 19091
$ws.Cells.Item(9, 16).Value = 1273

$ws.Cells.Item(10, 4).Value = 44449
$ws.Cells.Item(10, 10).Value = 220
$ws.Cells.Item(10, 11).Value = 22000
$ws.Cells.Item(10, 12).Value = 24000
$ws.Cells.Item(10, 13).Value = 23091
$ws.Cells.Item(10, 16).Value = 1539

$ws.Cells.Item(12, 4).Value = 44399
$ws.Cells.Item(12, 10).Value = 150
$ws.Cells.Item(12, 11).Value = 22000
$ws.Cells.Item(12, 12).Value = 22000
$ws.Cells.Item(12, 13).Value = 22000
$ws.Cells.Item(12, 16).Value = 1467

$ws.Cells.Item(13, 4).Value = 44391
$ws.Cells.Item(13, 10).Value = 160
$ws.Cells.Item(13, 11).Value = 20000
$ws.Cells.Item(13, 12).Value = 20000
$ws.Cells.Item(13, 13).Value = 20000
$ws.Cells.Item(13, 16).Value = 1333

$ws.Cells.Item(14, 4).Value = 44435
$ws.Cells.Item(14, 10).Value = 140
$ws.Cells.Item(14, 11).Value = 21000
$ws.Cells.Item(14, 12).Value = 23000
$ws.Cells.Item(14, 13).Value = 21714
$ws.Cells.Item(14, 16).Value = 1448

$ws.Cells.Item(15, 4).Value = 44476
$ws.Cells.Item(15, 10).Value = 220
$ws.Cells.Item(15, 11).Value = 20000
$ws.Cells.Item(15, 12).Value = 22000
$ws.Cells.Item(15, 13).Value = 20909
$ws.Cells.Item(15, 16).Value = 1394

$ws.Cells.Item(16, 4).Value = 44747
$ws.Cells.Item(16, 10).Value = 400
$ws.Cells.Item(16, 11).Value = 17000
$ws.Cells.Item(16, 12).Value = 19000
$ws.Cells.Item(16, 13).Value = 17850
$ws.Cells.Item(16, 16).Value = 1190

$ws.Cells.Item(17, 4).Value = 44365
$ws.Cells.Item(17, 10).Value = 580
$ws.Cells.Item(17, 11).Value = 20000
$ws.Cells.Item(17, 12).Value = 22000
$ws.Cells.Item(17, 13).Value = 21103
$ws.Cells.Item(17, 16).Value = 1407

$ws.Cells.Item(18, 4).Value = 44453
$ws.Cells.Item(18, 10).Value = 280
$ws.Cells.Item(18, 11).Value = 20000
$ws.Cells.Item(18, 12).Value = 22000
$ws.Cells.Item(18, 13).Value = 21286
$ws.Cells.Item(18, 16).Value = 1419

$ws.Cells.Item(19, 4).Value = 44398
$ws.Cells.Item(19, 10).Value = 130
$ws.Cells.Item(19, 11).Value = 20000
$ws.Cells.Item(19, 12).Value = 20000
$ws.Cells.Item(19, 13).Value = 20000
$ws.Cells.Item(19, 16).Value = 1333

$ws.Cells.Item(20, 4).Value = 44727
$ws.Cells.Item(20, 10).Value = 220
$ws.Cells.Item(20, 11).Value = 16000
$ws.Cells.Item(20, 12).Value = 18000
$ws.Cells.Item(20, 13).Value = 16909
$ws.Cells.Item(20, 16).Value = 1127

$ws.Cells.Item(21, 4).Value = 44748
$ws.Cells.Item(21, 10).Value = 200
$ws.Cells.Item(21, 11).Value = 16000
$ws.Cells.Item(21, 12).Value = 17000
$ws.Cells.Item(21, 13).Value = 16400
$ws.Cells.Item(21, 16).Value = 1093

$ws.Cells.Item(22, 4).Value = 44699
$ws.Cells.Item(22, 10).Value = 150
$ws.Cells.Item(22, 11).Value = 18000
$ws.Cells.Item(22, 12).Value = 20000
$ws.Cells.Item(22, 13).Value = 18667
$ws.Cells.Item(22, 16).Value = 1244

$ws.Cells.Item(23, 4).Value = 44396
$ws.Cells.Item(23, 10).Value = 130
$ws.Cells.Item(23, 11).Value = 22000
$ws.Cells.Item(23, 12).Value = 22000
$ws.Cells.Item(23, 13).Value = 22000
$ws.Cells.Item(23, 16).Value = 1467

$ws.Cells.Item(24, 4).Value = 44714
$ws.Cells.Item(24, 10).Value = 200
$ws.Cells.Item(24, 11).Value = 16000
$ws.Cells.Item(24, 12).Value = 17000
$ws.Cells.Item(24, 13).Value = 16400
$ws.Cells.Item(24, 16).Value = 1093

$ws.Cells.Item(25, 4).Value = 44406
$ws.Cells.Item(25, 10).Value = 400
$ws.Cells.Item(25, 11).Value = 20000
$ws.Cells.Item(25, 12).Value = 22000
$ws.Cells.Item(25, 13).Value = 20850
$ws.Cells.Item(25, 16).Value = 1390

